$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Row 1 header values (B1:E1)
$ws.Range("B1").Value = 15
$ws.Range("C1").Value = 16
$ws.Range("D1").Value = 15
$ws.Range("E1").Value = 16

# Row 2 ("CON") updated values (B2:E2)
$ws.Range("B2").Value = 108.29594798993637
$ws.Range("C2").Value = 105.61207287012184
$ws.Range("D2").Value = 107.06903482947129
$ws.Range("E2").Value = 107.76606483851549

# Row 3 ("STR") updated values (B3:E3)
$ws.Range("B3").Value = 107.03031794451725
$ws.Range("C3").Value = 104.88524901633632
$ws.Range("D3").Value = 105.8366287550395
$ws.Range("E3").Value = 105.58780782636271

# Update the selected range to match the new selection B1:E3
$ws.Range("B1:E3").Select()
